# Auto-generated Excel COM-interop script
# Updates computed market-price / profit columns (H-N) on each job sheet
# to match the latest scheduled data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10002
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10002
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10002
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10654
$ws.Range("H43").Value = 9499.75
$ws.Range("J43").Value = 5999.6665
$ws.Range("L43").Value = 5999.6665
$ws.Range("N43").Value = -6137.6665
$ws.Range("H59").Value = 4999.8
$ws.Range("I59").Value = 2499.5
$ws.Range("J59").Value = 6666.6665
$ws.Range("K59").Value = 7498.5
$ws.Range("L59").Value = 19999.9995
$ws.Range("M59").Value = -6941.5
$ws.Range("N59").Value = -21113.9995
$ws.Range("H62").Value = 9560.571
$ws.Range("I62").Value = 8996
$ws.Range("J62").Value = 10972
$ws.Range("K62").Value = 8996
$ws.Range("L62").Value = 10972
$ws.Range("M62").Value = -8372
$ws.Range("N62").Value = -12220
$ws.Range("H65").Value = 9560.571
$ws.Range("I65").Value = 8996
$ws.Range("J65").Value = 10972
$ws.Range("K65").Value = 44980
$ws.Range("L65").Value = 54860
$ws.Range("M65").Value = -41860
$ws.Range("N65").Value = -61100
$ws.Range("H86").Value = 3278.7778
$ws.Range("J86").Value = 3073.111
$ws.Range("L86").Value = 3073.111
$ws.Range("N86").Value = -5319.111
$ws.Range("H89").Value = 3278.7778
$ws.Range("J89").Value = 3073.111
$ws.Range("L89").Value = 15365.555
$ws.Range("N89").Value = -26597.555
$ws.Range("H101").Value = 1843.4286
$ws.Range("J101").Value = 382.5
$ws.Range("L101").Value = 1147.5
$ws.Range("N101").Value = -4391.5
$ws.Range("H113").Value = 2379.6875
$ws.Range("I113").Value = 2554.4285
$ws.Range("J113").Value = 2243.7778
$ws.Range("K113").Value = 2554.4285
$ws.Range("L113").Value = 2243.7778
$ws.Range("M113").Value = 699.5715
$ws.Range("N113").Value = -8751.7778
$ws.Range("H137").Value = 1667.1818
$ws.Range("I137").Value = 1234.1428
$ws.Range("K137").Value = 3702.4284
$ws.Range("M137").Value = -1152.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1881.1428
$ws.Range("I2").Value = 568.6
$ws.Range("J2").Value = 5162.5
$ws.Range("K2").Value = 568.6
$ws.Range("L2").Value = 5162.5
$ws.Range("M2").Value = -455.6
$ws.Range("N2").Value = -5388.5
$ws.Range("H32").Value = 3527.5
$ws.Range("I32").Value = 3527.5
$ws.Range("K32").Value = 3527.5
$ws.Range("M32").Value = -3240.5
$ws.Range("H61").Value = 2234.7144
$ws.Range("I61").Value = 1529.8
$ws.Range("K61").Value = 1529.8
$ws.Range("M61").Value = -1317.8
$ws.Range("H116").Value = 1881.1428
$ws.Range("I116").Value = 568.6
$ws.Range("J116").Value = 5162.5
$ws.Range("K116").Value = 568.6
$ws.Range("L116").Value = 5162.5
$ws.Range("M116").Value = 1725.4
$ws.Range("N116").Value = -9750.5
$ws.Range("H132").Value = 1283.3214
$ws.Range("I132").Value = 1274.3846
$ws.Range("J132").Value = 1399.5
$ws.Range("K132").Value = 3823.1538
$ws.Range("L132").Value = 4198.5
$ws.Range("M132").Value = -1293.1538
$ws.Range("N132").Value = -9258.5
$ws.Range("H136").Value = 2234.7144
$ws.Range("I136").Value = 1529.8
$ws.Range("K136").Value = 4589.4
$ws.Range("M136").Value = -2039.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1881.1428
$ws.Range("I3").Value = 568.6
$ws.Range("J3").Value = 5162.5
$ws.Range("K3").Value = 568.6
$ws.Range("L3").Value = 5162.5
$ws.Range("M3").Value = -454.6
$ws.Range("N3").Value = -5390.5
$ws.Range("H88").Value = 19329.666
$ws.Range("I88").Value = 8311
$ws.Range("J88").Value = 20707
$ws.Range("K88").Value = 8311
$ws.Range("L88").Value = 20707
$ws.Range("M88").Value = -7905
$ws.Range("N88").Value = -21519
$ws.Range("H91").Value = 19329.666
$ws.Range("I91").Value = 8311
$ws.Range("J91").Value = 20707
$ws.Range("K91").Value = 8311
$ws.Range("L91").Value = 20707
$ws.Range("M91").Value = -6907
$ws.Range("N91").Value = -23515

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1818.5
$ws.Range("I7").Value = 1427
$ws.Range("K7").Value = 1427
$ws.Range("M7").Value = -1314
$ws.Range("H22").Value = 1010.6667
$ws.Range("J22").Value = 1571.4286
$ws.Range("L22").Value = 1571.4286
$ws.Range("N22").Value = -2271.4286
$ws.Range("H132").Value = 2940.5789
$ws.Range("I132").Value = 3051.2942
$ws.Range("K132").Value = 9153.882599999999
$ws.Range("M132").Value = -6623.882599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 686.2414
$ws.Range("I2").Value = 773.25
$ws.Range("K2").Value = 773.25
$ws.Range("M2").Value = -660.25
$ws.Range("H26").Value = 39999
$ws.Range("J26").Value = 39999
$ws.Range("L26").Value = 39999
$ws.Range("N26").Value = -40559
$ws.Range("H50").Value = 39999
$ws.Range("J50").Value = 39999
$ws.Range("L50").Value = 39999
$ws.Range("N50").Value = -40995
$ws.Range("H126").Value = 4557.2666
$ws.Range("I126").Value = 4038.5557
$ws.Range("J126").Value = 5335.3335
$ws.Range("K126").Value = 12115.6671
$ws.Range("L126").Value = 16006.0005
$ws.Range("M126").Value = -9645.667099999999
$ws.Range("N126").Value = -20946.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1741.9375
$ws.Range("I7").Value = 1419.5714
$ws.Range("K7").Value = 1419.5714
$ws.Range("M7").Value = -1307.5714
$ws.Range("H30").Value = 437.2
$ws.Range("I30").Value = 437.2
$ws.Range("K30").Value = 437.2
$ws.Range("M30").Value = -329.2
$ws.Range("H46").Value = 2867
$ws.Range("I46").Value = 2322.3076
$ws.Range("K46").Value = 2322.3076
$ws.Range("M46").Value = -2134.3076
$ws.Range("H126").Value = 1741.9375
$ws.Range("I126").Value = 1419.5714
$ws.Range("K126").Value = 4258.7142
$ws.Range("M126").Value = -1788.7142
$ws.Range("H132").Value = 7459.727
$ws.Range("I132").Value = 6274.1816
$ws.Range("J132").Value = 8645.272000000001
$ws.Range("K132").Value = 18822.5448
$ws.Range("L132").Value = 25935.816
$ws.Range("M132").Value = -16292.5448
$ws.Range("N132").Value = -30995.816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 26000
$ws.Range("I58").Value = 45000
$ws.Range("J58").Value = 7000
$ws.Range("K58").Value = 45000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = -44692
$ws.Range("N58").Value = -7616
$ws.Range("H97").Value = 11857.333
$ws.Range("J97").Value = 11857.333
$ws.Range("L97").Value = 11857.333
$ws.Range("N97").Value = -13839.333
$ws.Range("H136").Value = 539.0303
$ws.Range("J136").Value = 498
$ws.Range("L136").Value = 1494
$ws.Range("N136").Value = -6594
